# RBA v2.4 - Alteração da pasta de origem
# Applies the text replacements captured by the commit diff:
#   - body (document.xml):  " 000110363365 - X " -> "  "
#   - body (document.xml):  "QWR"  -> "TERE"   (bold run)
#   - header (header1.xml): "QWER" -> "TRE"
#   - header (header1.xml): "QWR"  -> "TERE"
#   - header (header1.xml): "Qwer" -> "Tre"   (5 occurrences)
#   - header (header1.xml): "qwer" -> "tre"   (3 occurrences)

$d = $word.ActiveDocument

function Replace-AllText($range, [string]$find, [string]$replace, [bool]$wholeWord) {
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $range.Find.Execute($find, $true, $wholeWord, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# --- Body (document.xml) replacements ---
Replace-AllText $d.Content " 000110363365 - X " "  " $false
Replace-AllText $d.Content "QWR" "TERE" $true

# --- Header (header1.xml) replacements ---
$header = $d.Sections.First.Headers(1).Range

Replace-AllText $header "QWER" "TRE" $true
Replace-AllText $header "QWR" "TERE" $true
Replace-AllText $header "Qwer" "Tre" $true
Replace-AllText $header "qwer" "tre" $true
